{"js": "// The document contains the sentence:\n//   \"There are 7 available buildings, which are all upgradeable and all serve different needs:\"\n// and the edit changes the numeral \"7\" to the word \"eight\" (the building list was\n// extended to 8 entries - mine, quarry, sawmill, warehouse, market, farm, barracks, wall -\n// in an earlier commit; this commit just updates the count in the intro sentence).\n//\n// We scope the search to the full, unique sentence first so that we don't\n// accidentally touch any other occurrence of the digit \"7\" elsewhere in the\n// document (e.g. inside \"89476\" or \"1970s\"), then narrow the search to the\n// \"7\" inside that sentence and replace just that token with \"eight\".\n\nconst sentence = \"There are 7 available buildings, which are all upgradeable and all serve different needs:\";\n\nconst sentenceResults = context.document.body.search(sentence, { matchCase: true });\nsentenceResults.load(\"items\");\nawait context.sync();\n\nif (sentenceResults.items.length === 0) {\n  throw new Error(\"Could not find the target sentence in the document.\");\n}\n\nconst sentenceRange = sentenceResults.items[0];\n\nconst numberResults = sentenceRange.search(\"7\", { matchCase: true });\nnumberResults.load(\"items\");\nawait context.sync();\n\nif (numberResults.items.length === 0) {\n  throw new Error(\"Could not find the '7' to replace inside the target sentence.\");\n}\n\nnumberResults.items[0].insertText(\"eight\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The document contains the sentence:\n#   \"There are 7 available buildings, which are all upgradeable and all serve different needs:\"\n# The edit changes the numeral \"7\" to the word \"eight\" (the building list was\n# extended to 8 entries - mine, quarry, sawmill, warehouse, market, farm,\n# barracks, wall - in an earlier commit; this commit just updates the count\n# quoted in the intro sentence).\n#\n# We first locate the whole, unique sentence so we don't touch any other\n# occurrence of the digit \"7\" elsewhere in the document (e.g. inside \"89476\"\n# or \"1970s\"), then narrow down to a fresh Range scoped to just that sentence\n# and search again for \"7\" inside it, replacing only that token with \"eight\".\n\n$d = $word.ActiveDocument\n\n$sentenceRange = $d.Content\n$sentenceFind = $sentenceRange.Find\n$sentenceFind.Text = \"There are 7 available buildings, which are all upgradeable and all serve different needs:\"\n$sentenceFind.MatchCase = $true\n$sentenceFind.MatchWholeWord = $false\n$sentenceFound = $sentenceFind.Execute()\n\nif (-not $sentenceFound) {\n    throw \"Could not find the target sentence in the document.\"\n}\n\n# Work on a fresh Range scoped to the matched sentence so the subsequent\n# inner Find cannot escape it and match an unrelated \"7\" elsewhere.\n$numberRange = $d.Range($sentenceRange.Start, $sentenceRange.End)\n$numberFind = $numberRange.Find\n$numberFind.Text = \"7\"\n$numberFind.MatchCase = $true\n$numberFound = $numberFind.Execute()\n\nif (-not $numberFound) {\n    throw \"Could not find the '7' to replace inside the target sentence.\"\n}\n\n$numberRange.Text = \"eight\"\n"}
